# Season-record columns: add Wins / Losses / Ties next to the existing stats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header formatting used by the rest of row 1
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# PasteSpecial(xlPasteFormats) only touches formatting, but re-assert the text just in case
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-45) gets the team's season record: 79 wins, 83 losses, 0 ties
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 79  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 83  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
